$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.800.67"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.83"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.15"
$ws.Range("E5").Value = "  -4.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.671"
$ws.Range("E6").Value = "  -6.89%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.04"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.339"
$ws.Range("E9").Value = "  -6.00%  "
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0964"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "12.82"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.132.61"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.706"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.878.27"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.778.31"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.86"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("E19").Value = "  -5.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.70"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.50"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("E22").Value = "  -4.49%  "
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("E25").Value = "  -13.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.77"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.29"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.97"
$ws.Range("E28").Value = "  -4.21%  "
$ws.Range("E29").Value = "  -5.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.128.46"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.14"
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0568"
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.09"
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.823"
$ws.Range("E36").Value = "  -10.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  -5.10%  "
$ws.Range("E38").Value = "  -25.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "97.24"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.91"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0662"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0208"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.07"
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0830"
$ws.Range("E44").Value = "  +12.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.275.51"
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("E46").Value = "  -6.24%  "
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.89"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.24"
$ws.Range("E50").Value = "  -7.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.17"
$ws.Range("E51").Value = "  -6.72%  "
